# AD-EYE/TA/SimulinkConfig.xlsx - "New values of Goal and ip address"
#
# The BlockName column (A2:A28) previously stored values prefixed with
# "Mazda_RX8_Coupe_1/"; that per-vehicle prefix is dropped so only the
# bare block/topic name remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "Rviz"
$ws.Range("A3").Value  = "Map"
$ws.Range("A4").Value  = "Sensing"
$ws.Range("A5").Value  = "Localization"
$ws.Range("A6").Value  = "FakeLocalization"
$ws.Range("A7").Value  = "Detection"
$ws.Range("A8").Value  = "MissionPlanning"
$ws.Range("A9").Value  = "MotionPlanning"
$ws.Range("A10").Value = "Switch"
$ws.Range("A11").Value = "Ssmp"
$ws.Range("A12").Value = "GoalPoseX"
$ws.Range("A13").Value = "GoalPoseY"
$ws.Range("A14").Value = "GoalPoseZ"
$ws.Range("A15").Value = "GoalOrientX"
$ws.Range("A16").Value = "GoalOrientY"
$ws.Range("A17").Value = "GoalOrientZ"
$ws.Range("A18").Value = "GoalOrientW"
$ws.Range("A19").Value = "GoalTime"
$ws.Range("A20").Value = "GoalDistance"
$ws.Range("A21").Value = "GnssPoseSimulink"
$ws.Range("A22").Value = "PointsRawFloat32"
$ws.Range("A23").Value = "ImageRaw"
$ws.Range("A24").Value = "ClockFrequency"
$ws.Range("A25").Value = "SimulinkState"
$ws.Range("A26").Value = "CurrentVelocity"
$ws.Range("A27").Value = "PoseOtherCar"
$ws.Range("A28").Value = "CurrentPose"

# B15 / B18 used a redundant duplicate cell style (index 1) that is
# identical to the default style (index 0); re-apply a plain "General"
# number format so both cells fall back onto the shared default style.
$ws.Range("B15").NumberFormat = "General"
$ws.Range("B18").NumberFormat = "General"

# Column A widened slightly to better fit the shortened labels.
$ws.Columns.Item(1).ColumnWidth = 36.8333333333333

# Move the active cell back up to the top of the sheet.
[void]$ws.Range("A2:A28").Select()
[void]$ws.Range("A1").Select()
